# Uruguay Primera División workbook update
# - Updates odds data for rows 117-120 (existing matches re-scraped with fresh prices)
# - Finalises the placeholder match on row 193 (result now known) and
#   appends a brand-new upcoming fixture on row 194.
#
# Written as Excel COM-interop (PowerShell-style) against $excel.ActiveWorkbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $data) {
    foreach ($col in $data.Keys) {
        $ws.Range($col + $row).Value = $data[$col]
    }
}

# ---------------------------------------------------------------------------
# Rows 117-120: refreshed odds snapshots
# ---------------------------------------------------------------------------

Set-Row 117 @{
    "B" = 7013409
    "E" = "Nacional De Football"
    "F" = "Torque"
    "G" = 1
    "H" = 1
    "I" = "D"
    "J" = 1.666
    "K" = 3.9
    "L" = 4.5
    "M" = 1.615
    "N" = 4
    "O" = 4.75
    "P" = -0.75
    "Q" = 1.8
    "R" = 2.05
    "S" = 2.75
    "T" = 1.95
    "U" = 1.9
    "V" = -1
    "W" = 3
    "X" = -1
    "Y" = -1
    "Z" = 1.05
    "AA" = -1
    "AB" = 0.8999999999999999
}

Set-Row 118 @{
    "B" = 7013886
    "E" = "Racing Club de Montevideo"
    "F" = "Cerro"
    "G" = 0
    "H" = 1
    "I" = "A"
    "J" = 2.25
    "K" = 3.1
    "L" = 3.25
    "M" = 2.25
    "N" = 2.875
    "O" = 3.5
    "P" = -0.25
    "Q" = 1.95
    "R" = 1.9
    "S" = 2
    "T" = 1.925
    "U" = 1.925
    "V" = -1
    "W" = -1
    "X" = 2.5
    "Y" = -1
    "Z" = 0.8999999999999999
    "AA" = -1
    "AB" = 0.925
}

Set-Row 119 @{
    "B" = 7013885
    "E" = "La Luz"
    "F" = "Atletico Fenix Montevideo"
    "G" = 0
    "H" = 2
    "I" = "A"
    "J" = 3
    "K" = 3
    "L" = 2.4
    "M" = 2.9
    "N" = 2.75
    "O" = 2.6
    "P" = 0
    "Q" = 2.025
    "R" = 1.825
    "S" = 2
    "T" = 2.025
    "U" = 1.825
    "V" = -1
    "W" = -1
    "X" = 1.6
    "Y" = -1
    "Z" = 0.825
    "AA" = 0
    "AB" = 0
}

Set-Row 120 @{
    "B" = 7013702
    "E" = "Defensor Sporting"
    "F" = "Danubio"
    "G" = 0
    "H" = 2
    "I" = "A"
    "J" = 1.8
    "K" = 3.6
    "L" = 4.2
    "M" = 1.8
    "N" = 3.6
    "O" = 4.2
    "P" = -0.75
    "Q" = 2.05
    "R" = 1.8
    "S" = 2.25
    "T" = 1.85
    "U" = 2
    "V" = -1
    "W" = -1
    "X" = 3.2
    "Y" = -1
    "Z" = 0.8
    "AA" = -0.5
    "AB" = 0.5
}

# ---------------------------------------------------------------------------
# Row 193: the match has now been played - result + extra AH/OU columns land.
# B193 switches from the text placeholder id to the real numeric match id.
# ---------------------------------------------------------------------------

$ws.Range("B193").Value = 8110829

Set-Row 193 @{
    "G" = 2
    "H" = 1
    "I" = "H"
    "M" = 5.25
    "N" = 3.6
    "O" = 1.7
    "Q" = 1.975
    "R" = 1.875
    "T" = 2
    "U" = 1.85
    "V" = 4.25
    "W" = -1
    "X" = -1
    "Y" = 0.9750000000000001
    "Z" = -1
    "AA" = 1
    "AB" = -1
}

# ---------------------------------------------------------------------------
# Row 194: brand-new fixture appended at the bottom of the sheet.
# ---------------------------------------------------------------------------

Set-Row 194 @{
    "C" = "Uruguay Primera División"
    "E" = "Penarol"
    "F" = "CA River Plate"
    "J" = 1.333
    "K" = 5
    "L" = 9
    "M" = 1.4
    "N" = 4.75
    "O" = 7.5
    "P" = -1.25
    "Q" = 1.925
    "R" = 1.925
    "S" = 2.5
    "T" = 2
    "U" = 1.85
    "V" = 0
    "W" = 0
    "X" = 0
}

# A194 (bold/bordered "id row" style) and D194 (date style) copy their
# formatting from the row above so we reuse the existing style entries
# instead of fabricating new ones.
$ws.Range("A193").Copy() | Out-Null
$ws.Range("A194").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A194").Value = 192

$ws.Range("D193").Copy() | Out-Null
$ws.Range("D194").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D194").Value = 45411.83333333334

# B194 must stay a genuine *text* match-id (mirrors the old B193 placeholder)
# even though "8105807" parses as a number - force text via NumberFormat,
# then reset the format back to a clean/default style (copied from the
# already-unstyled, text-typed C193 cell) so no stray number format lingers
# on the cell itself.
$b194 = $ws.Range("B194")
$b194.NumberFormat = "@"
$b194.Value = "8105807"
$ws.Range("C193").Copy() | Out-Null
$b194.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$b194.Value = "8105807"

$excel.CutCopyMode = 0

Write-Host "Done."
